# daily auto push: 2026-01-30 14:01 UTC
#
# Inserts one new data row into Sheet1 at row 725 (between the existing
# "2026/01/30" block and the "2026/12/29" block), pushing every following
# row down by one. The new row records an additional "2026/01/30 20:00"
# sample (ranking 23). Dimension grows from A1:D766 to A1:D767.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 725; Excel shifts rows 725:766 down to 726:767.
$ws.Rows.Item(725).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/01/30"), not
# real date serials. Force a text number format before assigning the
# string so Excel's COM type-inference doesn't silently convert it to a
# date value, then drop back to the workbook's default "Normal" style so
# the cell doesn't carry a stray number-format override.
$ws.Range("A725").NumberFormat = "@"
$ws.Range("A725").Value = "2026/01/30"
$ws.Range("A725").Style = "Normal"

$ws.Range("B725").Value = "金"
$ws.Range("C725").Value = 20
$ws.Range("D725").Value = 23
